$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, report date range) ---
# Collapsing multi-run rich text into a single run is an accepted
# consequence of editing via the Characters/Value COM surface; the
# visible text is what changes (same font for every run already).
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Numeric cell updates across the precinct crime-stat table (rows 15-30) ---
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -41.666666666666
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = 31.081081081081
$ws.Range("L16").Value = 115.555555555556
$ws.Range("M16").Value = -17.796610169491
$ws.Range("N16").Value = -81.766917293233
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 11.764705882352
$ws.Range("I17").Value = 123
$ws.Range("J17").Value = 109
$ws.Range("K17").Value = 12.844036697247
$ws.Range("L17").Value = 26.80412371134
$ws.Range("M17").Value = 48.192771084337
$ws.Range("N17").Value = -37.244897959183
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -12.5
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = -26.190476190476
$ws.Range("M18").Value = 16.981132075471
$ws.Range("N18").Value = -85.238095238095
$ws.Range("C19").Value = 12
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 13.793103448275
$ws.Range("I19").Value = 272
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 13.333333333333
$ws.Range("L19").Value = 75.483870967741
$ws.Range("M19").Value = 48.63387978142
$ws.Range("N19").Value = -38.321995464852
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 76
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = 55.102040816326
$ws.Range("L20").Value = 145.161290322581
$ws.Range("M20").Value = 375
$ws.Range("N20").Value = -70.881226053639
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 1.298701298701
$ws.Range("I21").Value = 637
$ws.Range("J21").Value = 569
$ws.Range("K21").Value = 11.950790861159
$ws.Range("L21").Value = 71.698113207547
$ws.Range("M21").Value = 36.989247311828
$ws.Range("N21").Value = -66.098988823842
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 137.5
$ws.Range("I23").Value = 107
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = 30.487804878048
$ws.Range("L23").Value = 46.575342465753
$ws.Range("M23").Value = 81.355932203389
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 35
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = -5.405405405405
$ws.Range("I24").Value = 341
$ws.Range("J24").Value = 344
$ws.Range("K24").Value = -0.872093023255
$ws.Range("L24").Value = 24.908424908424
$ws.Range("M24").Value = -19.385342789598
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 19
$ws.Range("H25").Value = 35.714285714285
$ws.Range("I25").Value = 171
$ws.Range("J25").Value = 155
$ws.Range("K25").Value = 10.322580645161
$ws.Range("L25").Value = 33.59375
$ws.Range("M25").Value = -17.788461538461
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -55.555555555555
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = -12.121212121212
$ws.Range("L27").Value = 26.086956521739

# --- Cells that flip from a numeric 0 / a "N/A" percentage to the sheet's
#     text placeholders ("0" and "***.*") used when data is unavailable/zero.
#     Use PasteSpecial(Formats) from an existing placeholder cell so the
#     style index (General number format, right aligned) matches, then set
#     the literal text (prefixed with an apostrophe so it is stored as text,
#     not reinterpreted as a number).
$zeroSrc = $ws.Range("F28")
$naSrc = $ws.Range("N23")

$ws.Range("C15").Value = "'0"
$zeroSrc.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C18").Value = "'0"
$zeroSrc.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C22").Value = "'0"
$zeroSrc.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C26").Value = "'0"
$zeroSrc.Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$zeroSrc.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$naSrc.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G28").Value = "'0"
$zeroSrc.Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = "'***.*"
$naSrc.Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("G29").Value = "'0"
$zeroSrc.Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = "'***.*"
$naSrc.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("D30").Value = "'0"
$zeroSrc.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$naSrc.Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
